$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update stat values for teams that played on July 20 ---------------
# Columns: C = Los Angeles Lakers, E = Los Angeles Clippers, H = Utah Jazz,
# S = New Orleans Pelicans. Rows: 2 = FTA RATE, 3 = OPP EFG%,
# 4 = OPP FTA RATE, 5 = OPP TOV%, 6 = OPP OREB%.
$ws.Range("C2").Value = 0.45100000000000001
$ws.Range("E2").Value = 0.378
$ws.Range("H2").Value = 0.33300000000000002
$ws.Range("S2").Value = 0.191

$ws.Range("C3").Value = 54.1
$ws.Range("E3").Value = 45.7
$ws.Range("H3").Value = 48.4
$ws.Range("S3").Value = 48.8

$ws.Range("C4").Value = 0.378
$ws.Range("E4").Value = 0.45100000000000001
$ws.Range("H4").Value = 0.191
$ws.Range("S4").Value = 0.33300000000000002

$ws.Range("C5").Value = 21.6
$ws.Range("E5").Value = 15.5
$ws.Range("H5").Value = 20.8
$ws.Range("S5").Value = 19.600000000000001

$ws.Range("C6").Value = 18.600000000000001
$ws.Range("E6").Value = 25.5
$ws.Range("H6").Value = 36.799999999999997
$ws.Range("S6").Value = 34

# Mark every updated cell with the new "edited" font (Helvetica, dark grey)
$updated = "C2,E2,H2,S2,C3,E3,H3,S3,C4,E4,H4,S4,C5,E5,H5,S5,C6,E6,H6,S6"
foreach ($addr in $updated.Split(",")) {
    $ws.Range($addr).Font.Name = "Helvetica"
    $ws.Range($addr).Font.Size = 12
    $ws.Range($addr).Font.Color = 3355443
}

# --- 2. Drop the now-unused placeholder rows -------------------------------
# Rows 7-9 and 22-28 go back to being completely blank.
$ws.Range("A7:I9").Clear()
$ws.Range("A22:I28").Clear()

# --- 3. Trim the trailing placeholder rows down to columns A:D -------------
$ws.Range("I10:I13").Clear()
$ws.Range("E14:I17").Clear()
$ws.Range("E18:I21").Clear()

# Rows 10-17 col B:D lose their Helvetica formatting (back to default font)
$ws.Range("B10:D17").ClearFormats()

# Rows 18-21 col A picks up the plain (Helvetica) look already used by B:D
$ws.Range("B18:B21").Copy()
$ws.Range("A18:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Selection cursor ends on D9, matching the saved workbook ----------
$ws.Range("D9").Select()
